# Auto-generated Excel COM-interop script applying the weekly crime-data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -30
$ws.Range("L15").Value = 75
$ws.Range("N15").Value = -63.157894736842
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 25
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = -13.793103448275
$ws.Range("M16").Value = -55.357142857142
$ws.Range("N16").Value = -88.738738738738
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -48.387096774193
$ws.Range("I17").Value = 62
$ws.Range("J17").Value = 65
$ws.Range("K17").Value = -4.615384615384
$ws.Range("L17").Value = -4.615384615384
$ws.Range("M17").Value = 21.568627450980
$ws.Range("N17").Value = -49.593495934959
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -44.827586206896
$ws.Range("L18").Value = -48.387096774193
$ws.Range("M18").Value = -76.811594202898
$ws.Range("N18").Value = -90.303030303030
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 11
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 13.333333333333
$ws.Range("I19").Value = 75
$ws.Range("J19").Value = 73
$ws.Range("K19").Value = 2.739726027397
$ws.Range("L19").Value = 44.230769230769
$ws.Range("M19").Value = -33.035714285714
$ws.Range("N19").Value = -78.134110787172
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -41.379310344827
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -38.596491228070
$ws.Range("L20").Value = 29.629629629629
$ws.Range("M20").Value = -20.454545454545
$ws.Range("N20").Value = -88.817891373801
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -14.814814814814
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = -18.348623853211
$ws.Range("I21").Value = 220
$ws.Range("J21").Value = 257
$ws.Range("K21").Value = -14.396887159533
$ws.Range("L21").Value = 4.761904761904
$ws.Range("M21").Value = -34.131736526946
$ws.Range("N21").Value = -81.543624161073
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 50
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = 0
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 55.555555555555
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = -2.884615384615
$ws.Range("I24").Value = 239
$ws.Range("J24").Value = 265
$ws.Range("K24").Value = -9.811320754716
$ws.Range("L24").Value = 22.564102564102
$ws.Range("M24").Value = 32.044198895027
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -36.363636363636
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 109
$ws.Range("J25").Value = 86
$ws.Range("K25").Value = 26.744186046511
$ws.Range("L25").Value = 47.297297297297
$ws.Range("M25").Value = -26.845637583892
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 13
$ws.Range("K26").Value = -23.076923076923
$ws.Range("L26").Value = 25
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 22.222222222222
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("N28").Value = -84.848484848484
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -88.888888888888

# --- Cells changing between text ("N/A"-style shared string) and numeric ---
# Integer cells (style matches the #,##0 format used throughout the table)
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C18").Value = 2
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'

# Text placeholder cells ("0" / "***.*") replacing numeric values -
# set as text, then copy number-format/style from a stable same-style source cell
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122)
